$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: replace the "Jenn Cooper" contact with the new "Hannah" contact ---
# Clear the occupation/company/wechat_id columns entirely (they become blank for Hannah).
$ws.Range("D2:F2").ClearContents()

# Remove the existing hyperlinks first so they can be rebuilt in the order that
# matches the final file (Phil's C3 link becomes rId1, Hannah's C2 link becomes rId2).
$ws.Hyperlinks.Delete()

# Overwrite the first/last name and email for row 2.
$ws.Range("A2").Value = "Hannah"
$ws.Range("B2").Value = "Weier"
$ws.Range("C2").Value = "notnull@email.com"
$ws.Range("G2").Value = "empty fields except for fname lname email"

# --- Row 3 (Phil) keeps its data, just needs its hyperlink re-created ---
# (A3:E3 and G3 already hold the correct Phil Weier values.)

# Re-create the hyperlinks: Phil (C3) first, then Hannah (C2), so C3 -> rId1 and C2 -> rId2.
$ws.Hyperlinks.Add($ws.Range("C3"), "mailto:philweier@hotmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:notnull@email.com") | Out-Null

# Adding a hyperlink re-applies the built-in "Hyperlink" cell style but as a freshly
# duplicated style record; re-apply the named style explicitly so both email cells
# collapse back onto the workbook's single shared "Hyperlink" style.
$ws.Range("C2").Style = "Hyperlink"
$ws.Range("C3").Style = "Hyperlink"

# --- Update the selected cell to match the saved view state ---
$ws.Range("G2").Select() | Out-Null
